$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 441; this pushes rows 441..532 down to 442..533
$ws.Rows.Item(441).Insert()

# Populate the new row 441 with the new data record
$ws.Cells.Item(441, 1).Value = 4
$ws.Cells.Item(441, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(441, 3).Value = "Los Lagos"
$ws.Cells.Item(441, 4).Value = Get-Date -Year 2023 -Month 11 -Day 14 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(441, 5).Value = 10
$ws.Cells.Item(441, 6).Value = 100112043
$ws.Cells.Item(441, 7).Value = "Pepino ensalada"
$ws.Cells.Item(441, 8).Value = "Sin especificar"
$ws.Cells.Item(441, 9).Value = "Primera"
$ws.Cells.Item(441, 10).Value = 400
$ws.Cells.Item(441, 11).Value = 20000
$ws.Cells.Item(441, 12).Value = 20000
$ws.Cells.Item(441, 13).Value = 20000
$ws.Cells.Item(441, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(441, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(441, 16).Value = 333
$ws.Cells.Item(441, 17).Value = 60
$ws.Cells.Item(441, 18).Value = "Hortaliza"
